# Player-TEST.xlsx edit: clear the Team/Position/Nationality values that were
# accidentally left on row 2 (Tony Parker), then leave the selection where
# the user ended up navigating to (E9) after doing so.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").ClearContents()
$ws.Range("F2").ClearContents()
$ws.Range("J2").ClearContents()

$ws.Range("E9").Select()
